# adding TP & Labo Starter
# Replace the "*" placeholder markers with "p" across rows 11-17 on the
# route03 sheet (contiguous blocks of previously-"*" cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    "BV11:BY11",
    "AH12:AN12", "BV12:CA12",
    "AH13:AU13", "BT13:CA13",
    "AH14:AU14", "BT14:CA14",
    "AJ15:AU15", "BS15:CA15",
    "AJ16:AU16", "BS16:CA16",
    "AN17:AU17", "BS17:BX17"
)

foreach ($rangeAddress in $ranges) {
    $ws.Range($rangeAddress).Value = "p"
}
